# Applies the edit described by the diff:
#  1. Split "Dijk: ik somde miljarden op...  u bezuinigt kapot" around "u" with proofErr gramStart/gramEnd.
#  2. Split "Armoede kinderen / of / schrappen belasting aandelen" around "kinderen /" with proofErr gramStart/gramEnd.
#  3. Split ", helaas van cadeautjes" into ", helaas van cadeautje" + "s", and append a large
#     block of new meeting-notes paragraphs right after that paragraph.

$d = $word.ActiveDocument

function Replace-ParagraphXml($anchorText, $xmlPackage) {
    $rng = $d.Content
    $found = $rng.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "anchor not found: $anchorText"
    }
    $rng.Expand(4) | Out-Null   # wdParagraph - grow to the whole paragraph (incl. mark)
    $sub = $d.Range($rng.Start, $rng.End - 1)   # exclude the trailing paragraph mark
    $sub.InsertXML($xmlPackage)
}

$frag1 = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Dijk: ik somde miljarden op…  </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>u</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> bezuinigt kapot</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$frag2 = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r><w:tab/><w:t xml:space="preserve">Armoede </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>kinderen /</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> of / schrappen belasting aandelen</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$frag3 = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r><w:t xml:space="preserve">Ik kan een </w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>waslijst</w:t></w:r><w:r><w:t>, helaas van cadeautje</w:t></w:r><w:r><w:t>s</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t xml:space="preserve">Grote </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>vermogens /</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> grote bedrijven</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t xml:space="preserve">Potverdorie beter dan </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>yesilgos</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, geert is een </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>VVD leider</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p><w:r><w:tab/><w:t>Geert: Communistische teksten</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>Dijk: Ik word uitgedaagd (Dus wil toch terug)</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>Geert allemaal beloftes</w:t></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:tab/><w:t>Forse, u ben een loopjongen van het groot kapitaal (Beeldspraak)</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>Geert: geen antwoorden op vraag</w:t></w:r></w:p><w:p><w:r><w:tab/></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Timmermans, over </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>jimmy</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> dijk</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>Dijk wegzetten als communist is vreemd, u had zelfde beleid</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>U was uiterst kritisch op grootkapitaal</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Collegas</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> wegzetten</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:t>Timmermans herhaalt loopjongen VVD “quote van Dijk”</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Wilders: Dijk, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Shout</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Out van Wilders private </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>equity</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> in de zorg motie mee </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>akoord</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, geen winst in de zorg</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>15:20</w:t></w:r><w:r><w:tab/><w:t>botsen, voorbeelden verwijten</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>Zoeken naar punten om een over te worden</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>Onderhandelingen muurvast</w:t></w:r></w:p><w:p><w:r><w:t>Tijd om kant te kiezen</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>“Knieën versleten” kunnen niet wachten</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>Kant kiezen van mensen of grootkapitaal</w:t></w:r></w:p><w:p/><w:p/><w:p/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

Replace-ParagraphXml "Dijk: ik somde miljarden op" $frag1
Replace-ParagraphXml "Armoede kinderen / of / schrappen belasting aandelen" $frag2
Replace-ParagraphXml "Ik kan een " $frag3

Write-Output "edit applied"
